$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the first employee block (rows 16-21, DEIBIS BUELVAS VERGARA),
#    keeping only the second employee block (formerly rows 22-27, LIBARDO
#    MIGUEL BARRIOS SILVA) which shifts up to rows 16-21.
$ws.Rows("16:21").Delete()

# 2. The remaining block's "Periodo Mora" column (E) was in descending
#    order (1801..1708); the refreshed report lists it ascending
#    (1708..1801), so rewrite it explicitly.
$periods = @(1708, 1709, 1710, 1711, 1712, 1801)
for ($i = 0; $i -lt 6; $i++) {
    $r = 16 + $i
    $ws.Range("E$r").Value = $periods[$i]
}

# 3. Updated "Valor Mora" per period for the kept worker.
for ($r = 16; $r -le 21; $r++) {
    $ws.Range("G$r").Value = 781242
}

# 4. Updated aggregate "VALOR MORA" total.
$ws.Range("E11").Value = 177054

# 5. Updated "Cant. Trabajadores" (worker count) now that only one worker
#    remains in the statement.
$ws.Range("C13").Value = 1

# 6. The "Novedad de Retiro" / "Novedad de Ingreso" header columns were
#    swapped.
$ws.Range("H15").Value = "Novedad de Retiro"
$ws.Range("I15").Value = "Novedad de Ingreso"
